# "Working hour is updated"
#
# The sheet used to store the daily hour entries for Developer1/Developer2
# (columns B and C) as text values (e.g. the literal string "1.5") so that
# they would line up nicely with the right-aligned header. This pass turns
# those text entries into real numbers (right-aligned via direct formatting
# instead of a text number format) and appends one more day of hours
# (19 May 2021) for Developer1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# address -> new numeric value, for every cell that used to hold a
# shared-string text representation of a number.
$updates = [ordered]@{
    "C4"  = 1.5
    "B5"  = 1.5
    "C5"  = 1.5
    "B6"  = 0
    "C6"  = 2.5
    "B7"  = 2
    "C7"  = 2
    "B8"  = 0
    "C8"  = 1
    "B9"  = 3
    "C9"  = 2
    "B10" = 4
    "B11" = 1
    "B12" = 2.5
    "B13" = 4
    "B14" = 4
    "B15" = 8
    "B16" = 8
    "B17" = 5
}

$first = $true
foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)

    if ($first) {
        # Switch the cell away from the inherited "text" number format
        # (which is what forced the value to be stored as a string) back
        # to a plain General number format, keeping the right alignment.
        $cell.Style = "Normal"
        $cell.HorizontalAlignment = -4152   # xlRight
        $first = $false
    } else {
        # Reuse the exact formatting just created above for every other
        # cell instead of recreating it from scratch each time.
        $ws.Range("C4").Copy()
        $cell.PasteSpecial(-4122)           # xlPasteFormats
    }

    $cell.Value = $updates[$addr]
}

# New row of data: Developer1 logged 4 hours on 19 May 2021.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)        # xlPasteFormats (reuse date format)
$ws.Range("A18").Value = 44335

$ws.Range("C4").Copy()
$ws.Range("B18").PasteSpecial(-4122)        # xlPasteFormats (reuse number format)
$ws.Range("B18").Value = 4

$excel.CutCopyMode = $false

# Match the author's final selection.
$ws.Range("G16").Select()
